$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.834768176078796
$ws.Range("B1").Value = 1.823181986808777
$ws.Range("C1").Value = 2.076534986495972
$ws.Range("D1").Value = 3.606794357299805
$ws.Range("E1").Value = 3.897196769714355
